# Word had a plain coordinate-list layout with explicit en-US language
# formatting (pPr/rPr/lang) on every paragraph. The new layout adds a 10x5
# grid of (x, y, 0) coordinate strings - in a mix of ad-hoc textual styles
# (plain, comma separated, parenthesised, extra spaces...) - grouped by row
# (blank paragraph between rows), and drops the explicit language formatting
# that used to decorate every paragraph/run.
$d = $word.ActiveDocument

# Row 0 keeps its original first cell ("0, 0, 0"); rows are 10 coordinates
# wide, stepping by 10 on x, with y = 0, 10, 20, 30, 40.
$row0 = @('0, 0, 0', '(10 0 0)', '20,0,0', '(30, 0, 0)', '40 0 0', '50, 0 ,0', '(60 0 0)', '70,0 , 0', '80  ,0,0', '90 0 0')
$row1 = @('0 10 0', '10, 10, 0', '(20, 10, 0)', '30,10,0', '40 10 0', '(50 10 0)', '60 , 10 , 0', '(70,10,0)', '80 10 0', '90,10,0')
$row2 = @('(0, 20, 0)', '10 20 0', '20, 20 , 0', '(30 20 0)', '40 20 0', '50 , 20 , 0', '(60 20 0)', '70 20, 0', '80 20 0', '(90 20 0)')
$row3 = @('0 30 0', '(10,30,0)', '20,30,0', '30 30 0', '(40 30 0)', '50,30,0', '60,30 , 0', '70, 30 ,0', '80 30 0', '(90 30 0)')
$row4 = @('0 40 0', '10,40,0', '(20 40 0)', '30,40 ,0', '40 40 0', '50 40, 0', '(60 40 0)', '70,40,0', '80 40 0', '90,40,0')
$rows = @($row0, $row1, $row2, $row3, $row4)

# The cell that happens to land on the rendered page boundary keeps a
# <w:lastRenderedPageBreak/> marker in its run, ahead of the text.
$pageBreakRow = 2
$pageBreakCell = '(60 20 0)'

# Build the plain-paragraph body XML: one <w:p><w:r><w:t>...</w:t></w:r></w:p>
# per coordinate, with an empty <w:p/> between each row of 10 - and nothing
# else (no rPr/pPr/lang) so the paragraphs carry no direct formatting.
$bodyXml = New-Object System.Text.StringBuilder
for ($r = 0; $r -lt $rows.Count; $r++) {
    foreach ($cell in $rows[$r]) {
        $escaped = $cell.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')
        if ($r -eq $pageBreakRow -and $cell -eq $pageBreakCell) {
            [void]$bodyXml.Append("<w:p><w:r><w:lastRenderedPageBreak/><w:t>$escaped</w:t></w:r></w:p>")
        } else {
            [void]$bodyXml.Append("<w:p><w:r><w:t>$escaped</w:t></w:r></w:p>")
        }
    }
    if ($r -lt $rows.Count - 1) {
        [void]$bodyXml.Append('<w:p/>')
    }
}

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml.ToString() + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Drop every paragraph after the first one: the old document had 6
# paragraphs total, and everything from the 2nd one on goes away.
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $d.Paragraphs.Item($i).Range.Delete()
}

# Replace what remains (the original first paragraph, "0, 0, 0") by
# inserting the raw WordprocessingML package built above. Using InsertXML
# (rather than Range.Text / InsertParagraphAfter) means the new paragraphs
# get no inherited rPr/pPr - in particular no <w:lang val="en-US"/> - since
# the replacement content is taken verbatim from the supplied package
# instead of picking up formatting from the run/paragraph mark being
# replaced.
$d.Paragraphs.Item(1).Range.InsertXML($packageXml)
